# "Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta"
# The "Periodo Mora" column (E16:E51) listed billing periods newest-first
# (2109 down to 1810). This update re-sorts that list oldest-first
# (1810 up to 2109) as part of refreshing the EC (Estado de Cuenta) data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newPeriodOrder = @(
    "1810","1811","1812",
    "1901","1902","1903","1904","1905","1906","1907","1908","1909","1910","1911","1912",
    "2001","2002","2003","2004","2005","2006","2007","2008","2009","2010","2011","2012",
    "2101","2102","2103","2104","2105","2106","2107","2108","2109"
)

$startRow = 16
for ($i = 0; $i -lt $newPeriodOrder.Length; $i++) {
    $ws.Cells.Item($startRow + $i, 5).Value = $newPeriodOrder[$i]
}
